# Scheduled data refresh: re-pull Universalis market prices / profit
# calculations for each Leve-profit worksheet and write the refreshed
# columns H:N (currentAveragePrice*, LevePrice*, LeveProfit*) back in place.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1218.8462
$ws.Range("I40").Value = 1195.4166
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1195.4166
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1020.4166
$ws.Range("N40").Value = -1850

$ws.Range("H112").Value = 1153.4615
$ws.Range("I112").Value = 580.625
$ws.Range("J112").Value = 2070
$ws.Range("K112").Value = 1741.875
$ws.Range("L112").Value = 6210
$ws.Range("M112").Value = -633.875
$ws.Range("N112").Value = -8426

$ws.Range("H132").Value = 2469981.5
$ws.Range("I132").Value = 864.7727
$ws.Range("J132").Value = 111111110
$ws.Range("K132").Value = 2594.3181
$ws.Range("L132").Value = 333333330
$ws.Range("M132").Value = -64.31809999999996
$ws.Range("N132").Value = -333338390

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 605.5
$ws.Range("I15").Value = 605.5
$ws.Range("K15").Value = 605.5
$ws.Range("M15").Value = -255.5

$ws.Range("H32").Value = 7190700.5
$ws.Range("I32").Value = 1695282.9
$ws.Range("K32").Value = 1695282.9
$ws.Range("M32").Value = -1694995.9

$ws.Range("H61").Value = 4370082.5
$ws.Range("I61").Value = 2315820.5
$ws.Range("J61").Value = 11765426
$ws.Range("K61").Value = 2315820.5
$ws.Range("L61").Value = 11765426
$ws.Range("M61").Value = -2315608.5
$ws.Range("N61").Value = -11765850

$ws.Range("H74").Value = 48365100
$ws.Range("I74").Value = 37567030
$ws.Range("J74").Value = 106674670
$ws.Range("K74").Value = 37567030
$ws.Range("L74").Value = 106674670
$ws.Range("M74").Value = -37566156
$ws.Range("N74").Value = -106676418

$ws.Range("H77").Value = 48365100
$ws.Range("I77").Value = 37567030
$ws.Range("J77").Value = 106674670
$ws.Range("K77").Value = 187835150
$ws.Range("L77").Value = 533373350
$ws.Range("M77").Value = -187830782
$ws.Range("N77").Value = -533382086

$ws.Range("H110").Value = 1252.7778
$ws.Range("I110").Value = 1159.375
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1159.375
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 885.625
$ws.Range("N110").Value = -6090

$ws.Range("H132").Value = 14962459
$ws.Range("I132").Value = 15878700
$ws.Range("J132").Value = 11114249
$ws.Range("K132").Value = 47636100
$ws.Range("L132").Value = 33342747
$ws.Range("M132").Value = -47633570
$ws.Range("N132").Value = -33347807

$ws.Range("H136").Value = 4370082.5
$ws.Range("I136").Value = 2315820.5
$ws.Range("J136").Value = 11765426
$ws.Range("K136").Value = 6947461.5
$ws.Range("L136").Value = 35296278
$ws.Range("M136").Value = -6944911.5
$ws.Range("N136").Value = -35301378

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1940.58
$ws.Range("I86").Value = 1947.9694
$ws.Range("K86").Value = 1947.9694
$ws.Range("M86").Value = -824.9694

$ws.Range("H89").Value = 1940.58
$ws.Range("I89").Value = 1947.9694
$ws.Range("K89").Value = 9739.847
$ws.Range("M89").Value = -4123.847

$ws.Range("H107").Value = 550.6
$ws.Range("I107").Value = 562.5
$ws.Range("J107").Value = 384
$ws.Range("K107").Value = 562.5
$ws.Range("L107").Value = 384
$ws.Range("M107").Value = 1357.5
$ws.Range("N107").Value = -4224

$ws.Range("H134").Value = 22323560
$ws.Range("I134").Value = 29413428
$ws.Range("J134").Value = 5105312
$ws.Range("K134").Value = 88240284
$ws.Range("L134").Value = 15315936
$ws.Range("M134").Value = -88237749
$ws.Range("N134").Value = -15321006

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1968121.6
$ws.Range("I31").Value = 993327.6
$ws.Range("J31").Value = 5690062.5
$ws.Range("K31").Value = 993327.6
$ws.Range("L31").Value = 5690062.5
$ws.Range("M31").Value = -993032.6
$ws.Range("N31").Value = -5690652.5

$ws.Range("H34").Value = 1968121.6
$ws.Range("I34").Value = 993327.6
$ws.Range("J34").Value = 5690062.5
$ws.Range("K34").Value = 993327.6
$ws.Range("L34").Value = 5690062.5
$ws.Range("M34").Value = -993125.6
$ws.Range("N34").Value = -5690466.5

$ws.Range("H51").Value = 13000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 13000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 13000
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -14472

$ws.Range("H59").Value = 30108.857
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 32627
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 32627
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -34917

$ws.Range("H60").Value = 11098.182
$ws.Range("I60").Value = 11000
$ws.Range("J60").Value = 11102.857
$ws.Range("K60").Value = 11000
$ws.Range("L60").Value = 11102.857
$ws.Range("M60").Value = -10489
$ws.Range("N60").Value = -12124.857

$ws.Range("H61").Value = 13000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 13000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 13000
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -13696

$ws.Range("H132").Value = 3336574
$ws.Range("I132").Value = 5002271.5
$ws.Range("J132").Value = 5178.6
$ws.Range("K132").Value = 15006814.5
$ws.Range("L132").Value = 15535.8
$ws.Range("M132").Value = -15004284.5
$ws.Range("N132").Value = -20595.8

$ws.Range("H134").Value = 758327.1
$ws.Range("I134").Value = 3618.2896
$ws.Range("J134").Value = 2670256.2
$ws.Range("K134").Value = 10854.8688
$ws.Range("L134").Value = 8010768.600000001
$ws.Range("M134").Value = -8319.8688
$ws.Range("N134").Value = -8015838.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2394409.2
$ws.Range("I5").Value = 2404621.5
$ws.Range("J5").Value = 2382738
$ws.Range("K5").Value = 7213864.5
$ws.Range("L5").Value = 7148214
$ws.Range("M5").Value = -7213752.5
$ws.Range("N5").Value = -7148438

$ws.Range("H135").Value = 2394409.2
$ws.Range("I135").Value = 2404621.5
$ws.Range("J135").Value = 2382738
$ws.Range("K135").Value = 21641593.5
$ws.Range("L135").Value = 21444642
$ws.Range("M135").Value = -21639058.5
$ws.Range("N135").Value = -21449712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 14107558
$ws.Range("I132").Value = 17688238
$ws.Range("J132").Value = 9094608
$ws.Range("K132").Value = 53064714
$ws.Range("L132").Value = 27283824
$ws.Range("M132").Value = -53062184
$ws.Range("N132").Value = -27288884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1192.4412
$ws.Range("I40").Value = 849.7368
$ws.Range("J40").Value = 1626.5333
$ws.Range("K40").Value = 849.7368
$ws.Range("L40").Value = 1626.5333
$ws.Range("M40").Value = -713.7368
$ws.Range("N40").Value = -1898.5333

$ws.Range("H122").Value = 18441964
$ws.Range("I122").Value = 2130157.5
$ws.Range("J122").Value = 100001000
$ws.Range("K122").Value = 6390472.5
$ws.Range("L122").Value = 300003000
$ws.Range("M122").Value = -6388022.5
$ws.Range("N122").Value = -300007900

$ws.Range("H132").Value = 1552867.2
$ws.Range("I132").Value = 2021254
$ws.Range("J132").Value = 7190.6
$ws.Range("K132").Value = 6063762
$ws.Range("L132").Value = 21571.8
$ws.Range("M132").Value = -6061232
$ws.Range("N132").Value = -26631.8

$ws.Range("H136").Value = 1369026.1
$ws.Range("I136").Value = 1898464.1
$ws.Range("J136").Value = 1311.5834
$ws.Range("K136").Value = 5695392.300000001
$ws.Range("L136").Value = 3934.7502
$ws.Range("M136").Value = -5692842.300000001
$ws.Range("N136").Value = -9034.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1270.7916
$ws.Range("I122").Value = 1121.1666
$ws.Range("J122").Value = 1719.6666
$ws.Range("K122").Value = 3363.4998
$ws.Range("L122").Value = 5158.9998
$ws.Range("M122").Value = -913.4998000000001
$ws.Range("N122").Value = -10058.9998

$ws.Range("H132").Value = 1907268.1
$ws.Range("I132").Value = 1278817
$ws.Range("K132").Value = 3836451
$ws.Range("M132").Value = -3833921
